$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style (bold, centered, bordered) from an existing header cell (F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update existing values in row 2
$ws.Range("B2").Value = 0.03406936841555383
$ws.Range("C2").Value = 0.9996488542842338
$ws.Range("D2").Value = 0.1330350711633044
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5, n_estimators=50))])"

# New data cells
$ws.Range("G2").Value = 0.1194315095165318
$ws.Range("H2").Value = 0.989

$excel.CutCopyMode = 0
